$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16, 17, 18 undergo a cyclic rotation of their data:
#   new row16 <- old row17
#   new row17 <- old row18 (including the Activity/Aktivitet note columns K:N)
#   new row18 <- old row16 (loses the K:N columns, which only old row18 had)
#
# Columns A,B,D,E,F,G,H,Q,R are populated in all three rows already, so we can
# rotate their values directly. Columns K,L,M,N only exist (with content) on
# row 18 before the edit, and must move to row 17 afterwards.

$cols = @(1,2,4,5,6,7,8,17,18)  # A,B,D,E,F,G,H,Q,R

$row16vals = @{}
$row17vals = @{}
$row18vals = @{}
foreach ($c in $cols) {
    $row16vals[$c] = $ws.Cells.Item(16, $c).Value2
    $row17vals[$c] = $ws.Cells.Item(17, $c).Value2
    $row18vals[$c] = $ws.Cells.Item(18, $c).Value2
}
foreach ($c in $cols) {
    $ws.Cells.Item(16, $c).Value = $row17vals[$c]
    $ws.Cells.Item(17, $c).Value = $row18vals[$c]
    $ws.Cells.Item(18, $c).Value = $row16vals[$c]
}

# Move the K:N ("Ålder-Stadium", "Kön", "Aktivitet", "Metod") cells from row 18
# to row 17, then clear them from row 18.
$ws.Range("K18:N18").Copy($ws.Range("K17:N17"))
$ws.Range("K18:N18").ClearContents() | Out-Null
